# allroundersetlastseason_ipl.xlsx — "code and files for consistency and form"
#
# The sheet had several bowling-stat columns (Bowl_Inns..5W, i.e. P:Y) that
# were populated with a literal "-" placeholder text for all-rounders who
# never bowled. This pass replaces those placeholder dashes with the numeric
# convention used elsewhere in the sheet: 0 for most columns, and 100 for the
# Bowl_Ave / Bowl_SR columns (U, V / U, W depending on row) that divide by a
# zero wicket count.
#
# Two rows (2 and 44) had the whole P:Y block blank ("-"); many more rows
# only had U/W ("-") because V (Econ) was already a real number there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-StatCell {
    param(
        [string]$Addr,
        [double]$NewValue
    )
    $ws.Range($Addr).Value = $NewValue
}

# --- Row 2: full P:Y block was "-" -> zero out non-bowling-average/SR cols,
#     100 for the Bowl_Ave/Econ/Bowl_SR trio (U,V,W) ---
Set-StatCell "P2" 0
Set-StatCell "Q2" 0
Set-StatCell "R2" 0
Set-StatCell "S2" 0
Set-StatCell "T2" 0
# T2 also carried a one-off bold variant of the Tahoma font; normalise it to
# match the plain Tahoma style used by its row-mates (U2, V2, ...).
$ws.Range("T2").Font.Bold = $false
Set-StatCell "U2" 100
Set-StatCell "V2" 100
Set-StatCell "W2" 100
Set-StatCell "X2" 0
Set-StatCell "Y2" 0

# --- Rows where only Bowl_Ave (U) and Bowl_SR (W) were "-" ---
foreach ($row in 7, 10, 20, 22, 27, 28, 32, 33, 36, 38, 48) {
    Set-StatCell "U$row" 100
    Set-StatCell "W$row" 100
}

# --- Row 44: same full-block case as row 2 ---
Set-StatCell "P44" 0
Set-StatCell "Q44" 0
Set-StatCell "R44" 0
Set-StatCell "S44" 0
Set-StatCell "T44" 0
$ws.Range("T44").Font.Bold = $false
Set-StatCell "U44" 100
Set-StatCell "V44" 100
Set-StatCell "W44" 100
Set-StatCell "X44" 0
Set-StatCell "Y44" 0

# --- Restore the view/selection state recorded at last save ---
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("T52").Select()
